# Applies the commit "actualizando texto1 y agregando prueba1.xlsx":
# splits the bookmark out of the "AÑADIENDO SEGUNDA LINEA " paragraph and
# appends four new paragraphs (blank / text / blank / text), finishing with
# a new trailing blank paragraph that now owns the _GoBack bookmark.

$d = $word.ActiveDocument

# The existing paragraph that ends with "AÑADIENDO SEGUNDA LINEA " also
# carries the _GoBack bookmark. Remove it now so it can be re-created on
# the brand new trailing paragraph once that paragraph exists.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastOriginalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)

# Create five new paragraphs after it: blank, "#YA ESTA SUBIDO A GIT HUB",
# blank, "AGREGANDO NUEVAS LINEAS LUEGO DE HABERLO SUBIRLO A GITHUB", blank.
$d.Paragraphs.Add($lastOriginalParagraph.Range) | Out-Null
$d.Paragraphs.Add($d.Paragraphs.Item($d.Paragraphs.Count).Range) | Out-Null
$d.Paragraphs.Add($d.Paragraphs.Item($d.Paragraphs.Count).Range) | Out-Null
$d.Paragraphs.Add($d.Paragraphs.Item($d.Paragraphs.Count).Range) | Out-Null
$d.Paragraphs.Add($d.Paragraphs.Item($d.Paragraphs.Count).Range) | Out-Null

$paragraphCount = $d.Paragraphs.Count

# Fill in the text of the two non-blank new paragraphs (2nd and 4th of the
# five new ones), leaving the other three empty.
$d.Paragraphs.Item($paragraphCount - 3).Range.InsertBefore("#YA ESTA SUBIDO A GIT HUB")
$d.Paragraphs.Item($paragraphCount - 1).Range.InsertBefore("AGREGANDO NUEVAS LINEAS LUEGO DE HABERLO SUBIRLO A GITHUB")

# Re-create the _GoBack bookmark on the final (blank) paragraph.
$finalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $finalParagraph.Range)
